$d = $word.ActiveDocument

# --- 1. Fill in Programador # 6's name/carnet placeholder (first occurrence only) ---
$range = $d.Content
$range.Find.Execute("(Escribir nombre completo y número de carnet)", $true, $true, $false, $false, $false,
                     $true, 1, $false, "Keneth Abraham Velásquez Batz 2016-188", 1)

# --- 2. Fill in Programador # 6's assigned activity placeholder (first occurrence only) ---
$range2 = $d.Content
$range2.Find.Execute("(Escribir actividad asignada)", $true, $true, $false, $false, $false,
                      $true, 1, $false, "Modelo Empleado y EmpleadoDAO", 1)

# --- 3. Remove the _GoBack bookmark that trails "Modelo Venta y VentaDAO." ---
$word.Options.ShowHiddenBookmarks = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$d.Save()
